# Generate Report for handoff
# b.md has a new handoff pending: update Status to "Ready for handoff" on all
# sheets, and record the new Latest Handoff File / Latest Handoff Datetime for
# the b.md row on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: b.md row (row 3) status -> "Ready for handoff"
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# zh-cn sheet: b.md row (row 3)
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-20 03:24:20"

# Rebuild the hyperlinks for the zh-cn sheet so the C3 hyperlink's display
# text matches the new handoff file name (individual hyperlink edits are not
# supported in-place, so all hyperlinks are recreated with their original
# targets/text, except C3 which gets the new display text).
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/449c01a7e913b849b9e5e202ff0bbec1dc12f0a3/e2e/a.md.md", "", "", "a.md.md")
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5941e472fcf0c8a0430784272551b503732730ab/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/907bdc28482481d56e1d905c20edde4045bbe325/e2e/a.md.md", "", "", "a.md.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/15c13b947c27e7409342f14e6f87e7843d4311ed/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/449c01a7e913b849b9e5e202ff0bbec1dc12f0a3/e2e/b.md.md", "", "", "b.md.md")
$zhcn.Hyperlinks.Add($zhcn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5941e472fcf0c8a0430784272551b503732730ab/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/907bdc28482481d56e1d905c20edde4045bbe325/e2e/a.md.md", "", "", "a.md.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/15c13b947c27e7409342f14e6f87e7843d4311ed/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/449c01a7e913b849b9e5e202ff0bbec1dc12f0a3/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# de-de sheet: b.md row (row 3)
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$dede.Range("D3").Value = "2016-01-20 03:24:30"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/449c01a7e913b849b9e5e202ff0bbec1dc12f0a3/e2e/a.md.md", "", "", "a.md.md")
$dede.Hyperlinks.Add($dede.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dbf15a47f0c6276ce00f530400213b8b9cf497ea/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/f964b33276e6c8ee085c9c076801a03c29522418/e2e/a.md.md", "", "", "a.md.md")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/eed5454ef5fae82bc3fbade864006ff069234b31/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/449c01a7e913b849b9e5e202ff0bbec1dc12f0a3/e2e/b.md.md", "", "", "b.md.md")
$dede.Hyperlinks.Add($dede.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dbf15a47f0c6276ce00f530400213b8b9cf497ea/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/f964b33276e6c8ee085c9c076801a03c29522418/e2e/a.md.md", "", "", "a.md.md")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/eed5454ef5fae82bc3fbade864006ff069234b31/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/449c01a7e913b849b9e5e202ff0bbec1dc12f0a3/.localization-config", "", "", ".localization-config")
